$d = $word.ActiveDocument

# Locate the paragraph containing the "Requisitos" text followed by
# "LOQ4055: Quimica Inorgânica (Requisito fraco)", and remove the
# trailing empty paragraph plus the "Ver no Jupiter..." paragraph plus
# the "© 2020 ..." paragraph that follow it, while keeping the final
# empty paragraph and the page-break paragraph intact.

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Ver no Jupiter Salvar em pdf Salvar em docx") {
        # The empty paragraph right before this one is the first to remove.
        $startPara = $d.Paragraphs.Item($i - 1)
        $endPara = $p
    }
    if ($p.Range.Text -match [regex]::Escape("Contact: luizeleno@usp.br")) {
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $range = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $range.Delete()
}
